$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename two provinces to their abbreviated forms (values unchanged)
$ws.Range("A47").Value = "P.Nakhon S.Ayutthaya"
$ws.Range("A5").Value = "Bungkan"

# Re-sort the data range (header excluded) alphabetically by province name
$rng = $ws.Range("A2:B78")
$key1 = $ws.Range("A2:A78")
$rng.Sort($key1, 1)

# Widen column A and update the active selection / scroll position
$ws.Range("A1").ColumnWidth = 31.1
$ws.Range("E15").Select()
